$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy row style (A: bold/border/center, E: date format) to the two new rows ---
$ws.Range("A128").Copy()
$ws.Range("A129:A130").PasteSpecial(-4122)
$ws.Range("E128").Copy()
$ws.Range("E129:E130").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 127 updates
$ws.Cells.Item(127, 1).Value = 125
$ws.Cells.Item(127, 2).Value = 7011619
$ws.Cells.Item(127, 3).Value = "Azerbaijan Premier League"
$ws.Cells.Item(127, 4).Value = "Azerbaijan Premier League"
$ws.Cells.Item(127, 5).Value = 45360.33333333334
$ws.Cells.Item(127, 6).Value = "FK Kapaz"
$ws.Cells.Item(127, 7).Value = "Neftchi Baku"
$ws.Cells.Item(127, 8).Value = 3
$ws.Cells.Item(127, 9).Value = 3
$ws.Cells.Item(127, 10).Value = "D"
$ws.Cells.Item(127, 11).Value = 3.5
$ws.Cells.Item(127, 12).Value = 3.2
$ws.Cells.Item(127, 13).Value = 2
$ws.Cells.Item(127, 14).Value = 4
$ws.Cells.Item(127, 15).Value = 3.4
$ws.Cells.Item(127, 16).Value = 1.833
$ws.Cells.Item(127, 17).Value = 0.5
$ws.Cells.Item(127, 18).Value = 1.95
$ws.Cells.Item(127, 19).Value = 1.85
$ws.Cells.Item(127, 20).Value = 2.25
$ws.Cells.Item(127, 21).Value = 1.8
$ws.Cells.Item(127, 22).Value = 2
$ws.Cells.Item(127, 23).Value = -1
$ws.Cells.Item(127, 24).Value = 2.4
$ws.Cells.Item(127, 25).Value = -1
$ws.Cells.Item(127, 26).Value = 0.95
$ws.Cells.Item(127, 27).Value = -1
$ws.Cells.Item(127, 28).Value = 0.8
$ws.Cells.Item(127, 29).Value = -1

# Row 128 updates
$ws.Cells.Item(128, 1).Value = 126
$ws.Cells.Item(128, 2).Value = 7011621
$ws.Cells.Item(128, 3).Value = "Azerbaijan Premier League"
$ws.Cells.Item(128, 4).Value = "Azerbaijan Premier League"
$ws.Cells.Item(128, 5).Value = 45360.4375
$ws.Cells.Item(128, 6).Value = "Sabah"
$ws.Cells.Item(128, 7).Value = "FK Gabala"
$ws.Cells.Item(128, 8).Value = 2
$ws.Cells.Item(128, 9).Value = 1
$ws.Cells.Item(128, 10).Value = "H"
$ws.Cells.Item(128, 11).Value = 1.4
$ws.Cells.Item(128, 12).Value = 4.333
$ws.Cells.Item(128, 13).Value = 6
$ws.Cells.Item(128, 14).Value = 1.533
$ws.Cells.Item(128, 15).Value = 4
$ws.Cells.Item(128, 16).Value = 4.75
$ws.Cells.Item(128, 17).Value = -1
$ws.Cells.Item(128, 18).Value = 1.925
$ws.Cells.Item(128, 19).Value = 1.875
$ws.Cells.Item(128, 20).Value = 2.5
$ws.Cells.Item(128, 21).Value = 1.825
$ws.Cells.Item(128, 22).Value = 1.975
$ws.Cells.Item(128, 23).Value = 0.5329999999999999
$ws.Cells.Item(128, 24).Value = -1
$ws.Cells.Item(128, 25).Value = -1
$ws.Cells.Item(128, 26).Value = 0
$ws.Cells.Item(128, 27).Value = 0
$ws.Cells.Item(128, 28).Value = 0.825
$ws.Cells.Item(128, 29).Value = -1

# Row 129 (new)
$ws.Cells.Item(129, 1).Value = 127
$ws.Cells.Item(129, 2).Value = 7011617
$ws.Cells.Item(129, 3).Value = "Azerbaijan Premier League"
$ws.Cells.Item(129, 4).Value = "Azerbaijan Premier League"
$ws.Cells.Item(129, 5).Value = 45361.35416666666
$ws.Cells.Item(129, 6).Value = "PFK Turan Tovuz"
$ws.Cells.Item(129, 7).Value = "Araz FK"
$ws.Cells.Item(129, 8).Value = 2
$ws.Cells.Item(129, 9).Value = 0
$ws.Cells.Item(129, 10).Value = "H"
$ws.Cells.Item(129, 11).Value = 2.3
$ws.Cells.Item(129, 12).Value = 3.1
$ws.Cells.Item(129, 13).Value = 2.8
$ws.Cells.Item(129, 14).Value = 2.05
$ws.Cells.Item(129, 15).Value = 3.25
$ws.Cells.Item(129, 16).Value = 3.2
$ws.Cells.Item(129, 17).Value = -0.25
$ws.Cells.Item(129, 18).Value = 1.825
$ws.Cells.Item(129, 19).Value = 1.975
$ws.Cells.Item(129, 20).Value = 2.25
$ws.Cells.Item(129, 21).Value = 1.8
$ws.Cells.Item(129, 22).Value = 2
$ws.Cells.Item(129, 23).Value = 1.05
$ws.Cells.Item(129, 24).Value = -1
$ws.Cells.Item(129, 25).Value = -1
$ws.Cells.Item(129, 26).Value = 0.825
$ws.Cells.Item(129, 27).Value = -1
$ws.Cells.Item(129, 28).Value = -0.5
$ws.Cells.Item(129, 29).Value = 0.5

# Row 130 (new)
$ws.Cells.Item(130, 1).Value = 128
$ws.Cells.Item(130, 2).Value = 7011618
$ws.Cells.Item(130, 3).Value = "Azerbaijan Premier League"
$ws.Cells.Item(130, 4).Value = "Azerbaijan Premier League"
$ws.Cells.Item(130, 5).Value = 45361.45833333334
$ws.Cells.Item(130, 6).Value = "FK Qarabag"
$ws.Cells.Item(130, 7).Value = "Zira IK"
$ws.Cells.Item(130, 8).Value = 3
$ws.Cells.Item(130, 9).Value = 1
$ws.Cells.Item(130, 10).Value = "H"
$ws.Cells.Item(130, 11).Value = 1.4
$ws.Cells.Item(130, 12).Value = 4
$ws.Cells.Item(130, 13).Value = 7
$ws.Cells.Item(130, 14).Value = 1.55
$ws.Cells.Item(130, 15).Value = 3.6
$ws.Cells.Item(130, 16).Value = 5.25
$ws.Cells.Item(130, 17).Value = -1
$ws.Cells.Item(130, 18).Value = 1.95
$ws.Cells.Item(130, 19).Value = 1.85
$ws.Cells.Item(130, 20).Value = 2.5
$ws.Cells.Item(130, 21).Value = 2
$ws.Cells.Item(130, 22).Value = 1.8
$ws.Cells.Item(130, 23).Value = 0.55
$ws.Cells.Item(130, 24).Value = -1
$ws.Cells.Item(130, 25).Value = -1
$ws.Cells.Item(130, 26).Value = 0.95
$ws.Cells.Item(130, 27).Value = -1
$ws.Cells.Item(130, 28).Value = 1
$ws.Cells.Item(130, 29).Value = -1
